$p = $ppt.ActivePresentation

$oldStyleId = "{1229C4F3-53D0-490D-9692-5A7939145234}"
$newStyleId = "{0EFFEE52-C7BA-449E-925E-BB50F22524E2}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $shp = $s.Shapes.Item($shi)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            $tbl.ApplyStyle($newStyleId)
        }
    }
}
